$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7447250485420227
$ws.Range("B1").Value = 2.8744957447052
$ws.Range("C1").Value = 3.113898515701294
$ws.Range("D1").Value = 3.719011068344116
$ws.Range("E1").Value = 1.301419138908386
